$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 1834
    $ws.Range("F6").Value = 1066
    $ws.Range("F8").Value = 5932
}
